$d = $word.ActiveDocument

# Step 1: Replace "normal pattern" with "standard reduction pattern"
$d.Content.Find.Execute("normal pattern", $true, $false, $false, $false, $false, $true, 1, $false, "standard reduction pattern", 2)
